# Update the "Förändrad" (changed) date column (C) for rows 2-18
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C18").Value = 45233
